# Updated cryptos list on Tue Jun 25 08:55:01 UTC 2024 with GitHub Actions
# Applies refreshed price/volume figures (and a Mantle/ONDO row swap) to Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '60.781.25'
$ws.Cells.Item(2, 5).Value = '  -3.37%  '

$ws.Cells.Item(3, 4).Value = '3.350.83'
$ws.Cells.Item(3, 5).Value = '  -1.16%  '

$ws.Cells.Item(4, 5).Value = '  +0.02%  '

$ws.Cells.Item(5, 4).Value = '''569.50'
$ws.Cells.Item(5, 5).Value = '  -0.99%  '

$ws.Cells.Item(6, 4).Value = '''133.78'
$ws.Cells.Item(6, 5).Value = '  +5.92%  '

$ws.Cells.Item(7, 5).Value = '  -0.04%  '

$ws.Cells.Item(8, 4).Value = '3.352.44'
$ws.Cells.Item(8, 5).Value = '  -1.16%  '

$ws.Cells.Item(9, 4).Value = '''0.476'
$ws.Cells.Item(9, 5).Value = '  -0.30%  '

$ws.Cells.Item(10, 5).Value = '  +3.39%  '

$ws.Cells.Item(11, 5).Value = '  +2.17%  '

$ws.Cells.Item(12, 5).Value = '  +2.69%  '

$ws.Cells.Item(13, 4).Value = '3.921.66'
$ws.Cells.Item(13, 5).Value = '  -1.11%  '

$ws.Cells.Item(14, 5).Value = '  +1.77%  '

$ws.Cells.Item(15, 5).Value = '  -0.24%  '

$ws.Cells.Item(16, 4).Value = '3.348.48'
$ws.Cells.Item(16, 5).Value = '  -1.25%  '

$ws.Cells.Item(17, 4).Value = '''24.92'
$ws.Cells.Item(17, 5).Value = '  +0.69%  '

$ws.Cells.Item(18, 4).Value = '60.850.73'

$ws.Cells.Item(19, 4).Value = '''13.86'
$ws.Cells.Item(19, 5).Value = '  +5.39%  '

$ws.Cells.Item(20, 5).Value = '  +1.91%  '

$ws.Cells.Item(21, 4).Value = '''9.33'
$ws.Cells.Item(21, 5).Value = '  +0.60%  '

$ws.Cells.Item(22, 4).Value = '''372.27'
$ws.Cells.Item(22, 5).Value = '  +0.07%  '

$ws.Cells.Item(23, 5).Value = '  +2.60%  '

$ws.Cells.Item(24, 4).Value = '3.484.32'
$ws.Cells.Item(24, 5).Value = '  -1.13%  '

$ws.Cells.Item(25, 5).Value = '  +0.08%  '

$ws.Cells.Item(26, 4).Value = '''70.36'
$ws.Cells.Item(26, 5).Value = '  -1.92%  '

$ws.Cells.Item(27, 5).Value = '  +9.27%  '

$ws.Cells.Item(28, 5).Value = '  +16.08%  '

$ws.Cells.Item(29, 4).Value = '''7.64'
$ws.Cells.Item(29, 5).Value = '  +8.38%  '

$ws.Cells.Item(30, 5).Value = '  +0.33%  '

$ws.Cells.Item(31, 4).Value = '''8.00'
$ws.Cells.Item(31, 5).Value = '  +1.81%  '

$ws.Cells.Item(32, 5).Value = '  +0.37%  '

$ws.Cells.Item(33, 5).Value = '  +2.78%  '

$ws.Cells.Item(34, 5).Value = '  -0.06%  '

$ws.Cells.Item(35, 4).Value = '3.380.60'
$ws.Cells.Item(35, 5).Value = '  -1.09%  '

$ws.Cells.Item(36, 4).Value = '''23.32'
$ws.Cells.Item(36, 5).Value = '  +2.58%  '

$ws.Cells.Item(37, 5).Value = '  +1.75%  '

$ws.Cells.Item(38, 5).Value = '  +2.70%  '

$ws.Cells.Item(39, 5).Value = '  +2.90%  '

$ws.Cells.Item(40, 4).Value = '''162.48'
$ws.Cells.Item(40, 5).Value = '  -2.34%  '

$ws.Cells.Item(41, 4).Value = '''0.0782'
$ws.Cells.Item(41, 5).Value = '  +3.36%  '

$ws.Cells.Item(42, 4).Value = '''1.00'
$ws.Cells.Item(42, 5).Value = '  +0.06%  '

$ws.Cells.Item(43, 4).Value = '''41.21'
$ws.Cells.Item(43, 5).Value = '  -1.63%  '

$ws.Cells.Item(44, 5).Value = '  +2.97%  '

$ws.Cells.Item(45, 2).Value = 'ONDO'
$ws.Cells.Item(45, 3).Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Cells.Item(45, 4).Value = '''1.20'
$ws.Cells.Item(45, 5).Value = '  +8.23%  '

$ws.Cells.Item(46, 2).Value = 'Mantle'
$ws.Cells.Item(46, 3).Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Cells.Item(46, 4).Value = '''0.755'
$ws.Cells.Item(46, 5).Value = '  -1.23%  '

$ws.Cells.Item(47, 5).Value = '  +2.14%  '

$ws.Cells.Item(48, 5).Value = '  +4.84%  '

$ws.Cells.Item(49, 4).Value = '''22.63'
$ws.Cells.Item(49, 5).Value = '  +0.62%  '

$ws.Cells.Item(50, 4).Value = '''23.12'
$ws.Cells.Item(50, 5).Value = '  +11.96%  '

$ws.Cells.Item(51, 5).Value = '  +12.23%  '
